$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 185
$ws.Range("I12").Value = 137.5
$ws.Range("K12").Value = 137.5
$ws.Range("M12").Value = 32.5
$ws.Range("H80").Value = 1319.2
$ws.Range("I80").Value = 1459
$ws.Range("K80").Value = 4377
$ws.Range("M80").Value = -3379
$ws.Range("H83").Value = 1319.2
$ws.Range("I83").Value = 1459
$ws.Range("K83").Value = 13131
$ws.Range("M83").Value = -8139
$ws.Range("H92").Value = 265.35715
$ws.Range("I92").Value = 275
$ws.Range("K92").Value = 275
$ws.Range("M92").Value = 973
$ws.Range("H111").Value = 850
$ws.Range("I111").Value = 450
$ws.Range("K111").Value = 1350
$ws.Range("M111").Value = 1717
$ws.Range("H113").Value = 1893
$ws.Range("I113").Value = 1841.25
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 1841.25
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 1412.75
$ws.Range("N113").Value = -8608
$ws.Range("H141").Value = 5795
$ws.Range("I141").Value = 5795
$ws.Range("K141").Value = 17385
$ws.Range("M141").Value = -12205

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 900
$ws.Range("I30").Value = 900
$ws.Range("K30").Value = 900
$ws.Range("M30").Value = -750
$ws.Range("H45").Value = 1900
$ws.Range("J45").Value = 2000
$ws.Range("L45").Value = 2000
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("H122").Value = 6653.8
$ws.Range("I122").Value = 5923
$ws.Range("J122").Value = 7750
$ws.Range("K122").Value = 17769
$ws.Range("L122").Value = 23250
$ws.Range("M122").Value = -15319
$ws.Range("N122").Value = -28150
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 50
$ws.Range("I80").Value = 50
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 50
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 948
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 50
$ws.Range("I83").Value = 50
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 250
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 4742
$ws.Range("N83").ClearContents()
$ws.Range("H86").Value = 7000
$ws.Range("I86").Value = 7000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 7000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -5877
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 7000
$ws.Range("I89").Value = 7000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 35000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -29384
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 1767.9445
$ws.Range("I107").Value = 1723.5333
$ws.Range("K107").Value = 1723.5333
$ws.Range("M107").Value = 196.4666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 38333
$ws.Range("J59").Value = 45000
$ws.Range("L59").Value = 45000
$ws.Range("N59").Value = -47290
$ws.Range("H99").Value = 600919.8
$ws.Range("I99").Value = 1399.5
$ws.Range("J99").Value = 1000600
$ws.Range("K99").Value = 1399.5
$ws.Range("L99").Value = 1000600
$ws.Range("M99").Value = 98.5
$ws.Range("N99").Value = -1003596
$ws.Range("H103").Value = 20262
$ws.Range("I103").Value = 20262
$ws.Range("K103").Value = 20262
$ws.Range("M103").Value = -19090
$ws.Range("H126").Value = 600919.8
$ws.Range("I126").Value = 1399.5
$ws.Range("J126").Value = 1000600
$ws.Range("K126").Value = 4198.5
$ws.Range("L126").Value = 3001800
$ws.Range("M126").Value = -1728.5
$ws.Range("N126").Value = -3006740
$ws.Range("H134").Value = 1826.0834
$ws.Range("I134").Value = 1782
$ws.Range("K134").Value = 5346
$ws.Range("M134").Value = -2811

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 132
$ws.Range("I7").Value = 77
$ws.Range("J7").Value = 177.83333
$ws.Range("K7").Value = 231
$ws.Range("L7").Value = 533.49999
$ws.Range("M7").Value = -119
$ws.Range("N7").Value = -757.49999
$ws.Range("H117").Value = 1239.8
$ws.Range("I117").Value = 800
$ws.Range("J117").Value = 1349.75
$ws.Range("K117").Value = 2400
$ws.Range("L117").Value = 4049.25
$ws.Range("M117").Value = 1042
$ws.Range("N117").Value = -10933.25
$ws.Range("H121").Value = 2000
$ws.Range("J121").Value = 2000
$ws.Range("L121").Value = 6000
$ws.Range("N121").Value = -8620
$ws.Range("H129").Value = 2344.4285
$ws.Range("I129").Value = 344.5
$ws.Range("K129").Value = 1033.5
$ws.Range("M129").Value = 3966.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H107").Value = 1578.6666
$ws.Range("J107").Value = 2121.6
$ws.Range("L107").Value = 2121.6
$ws.Range("N107").Value = -5961.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4816.5454
$ws.Range("J46").Value = 4816.5454
$ws.Range("L46").Value = 4816.5454
$ws.Range("N46").Value = -5192.5454
$ws.Range("H132").Value = 2481.5
$ws.Range("I132").Value = 2481.5
$ws.Range("K132").Value = 7444.5
$ws.Range("M132").Value = -4914.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 78000
$ws.Range("J16").Value = 78000
$ws.Range("L16").Value = 78000
$ws.Range("N16").Value = -78584
$ws.Range("H62").Value = 3525
$ws.Range("I62").Value = 3525
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3525
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2901
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3525
$ws.Range("I65").Value = 3525
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 17625
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -14505
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H107").Value = 694.75
$ws.Range("I107").Value = 694.75
$ws.Range("K107").Value = 2084.25
$ws.Range("M107").Value = -164.25
$ws.Range("H113").Value = 5917.7
$ws.Range("I113").Value = 397.125
$ws.Range("K113").Value = 1191.375
$ws.Range("M113").Value = 978.625
$ws.Range("H122").Value = 1599.75
$ws.Range("I122").Value = 1299.6666
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3898.9998
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1448.9998
$ws.Range("N122").Value = -12400
